# Insert a new data row above row 216 (shifting all subsequent rows down by one)
# and populate it with a copy of the row that lands at 217 after the shift,
# except for the Fecha (D) and Volumen (J) values which take new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("216:216").Insert()

# After the insert, the data that used to live in row 216 now lives in row 217.
# Copy that row's values into the freshly inserted row 216 for every column
# except D (Fecha) and J (Volumen), which get their own new values below.
for ($col = 1; $col -le 18; $col++) {
    if ($col -eq 4 -or $col -eq 10) { continue }
    $ws.Cells.Item(216, $col).Value = $ws.Cells.Item(217, $col).Value2
}

$ws.Cells.Item(216, 4).Value = 44606
$ws.Cells.Item(216, 10).Value = 100
